$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = "S2 0"; "C2" = "S2 1"; "D2" = "S3 1"; "E2" = "S3 0"; "F2" = "S2 0"; "G2" = "S2 1"; "H2" = "S3 1"; "I2" = "S3 0";
    "B3" = "S3 1"; "D3" = "S0 1"; "E3" = "S0 0"; "F3" = "S3 0"; "H3" = "S0 0"; "I3" = "S0 1";
    "B4" = "S0 1"; "C4" = "S2 0"; "D4" = "S0 1"; "E4" = "S2 0"; "F4" = "S0 0"; "G4" = "S1 0"; "I4" = "S1 0";
    "C5" = "S2 1"; "E5" = "S2 1"; "F5" = "S1 0"; "G5" = "S2 1"; "H5" = "S1 0"; "I5" = "S2 1"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws.Range("H4").Select()
